# DT added label option
# Refresh the OLS regression result blobs (re-run timestamp) on each of the
# three "backward elimination" summary sheets: the Date/Time line inside the
# big text dump in column B, row 2 is updated to reflect the new run.

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 05 Jan 2020"
$newDate = "Wed, 08 Jan 2020"
$oldTime = "21:22:19"

# Sheet "5" (F-statistic 7.970) and sheet "4" (F-statistic 9.951) were
# re-run at 19:07:24; sheet "3" (F-statistic 17.43) finished a second later,
# at 19:07:25.
$newTimeBySheet = @{
    "5" = "19:07:24"
    "4" = "19:07:24"
    "3" = "19:07:25"
}

foreach ($ws in $wb.Worksheets) {
    $newTime = $newTimeBySheet[$ws.Name]
    if (-not $newTime) { continue }

    $cell = $ws.Range("B2")
    $text = $cell.Value2
    if ($null -eq $text) { continue }

    $text = $text.Replace($oldDate, $newDate)
    $text = $text.Replace($oldTime, $newTime)

    $cell.Value2 = $text
}
